$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Colors (BGR integers for Interior.Color)
#   light green FF92D050 (existing highlight) -> already used by E3
#   dark  green FF00B050 (new highlight)      -> used by most formula cells
# ---------------------------------------------------------------------------
$lightGreen = 0x50 * 65536 + 0xD0 * 256 + 0x92   # FF92D050
$darkGreen  = 0x50 * 65536 + 0xB0 * 256 + 0x00   # FF00B050

# ---------------------------------------------------------------------------
# 1. Clear out the region that is being completely re-laid-out (old rows
#    1-15, columns A and B, plus G25) so stale cells don't linger.
# ---------------------------------------------------------------------------
$ws.Range("A1:B16").Clear()
$ws.Range("E1:E3").Clear()
$ws.Range("G24:G27").ClearContents()

# ---------------------------------------------------------------------------
# 2. Re-enter the formulas in their new positions (exact source text - the
#    engine re-normalises redundant unary +/- and ^ grouping on save, same
#    as it would for any workbook it round-trips).
# ---------------------------------------------------------------------------
$ws.Range("A1").Formula  = "=1+4*2-4^2+2^2^2*2+3-1"
$ws.Range("A2").Formula  = "=((1+3-2/2)*(2^2^2)^2+1)*2^2-1"
$ws.Range("A3").Formula  = "=1+4*2-4*2+2*2*2*2+3-1"
$ws.Range("A4").Formula  = "=1+4*--2-4*+-+-2+2*--2*2*2+3-++-+--1"
$ws.Range("A5").Formula  = "=3+(2+(3*2)*2+(2^2)--4)"
$ws.Range("A6").Formula  = "=1*(2^2^2)^2+1"
$ws.Range("A7").Formula  = "=(1)*(2^2^2)^2+1"
$ws.Range("A8").Formula  = "=((1+1)*(2^2^2)^2+(1--+((1+1)*4)))"
$ws.Range("A9").Formula  = "=(1+1)*(3^3^2)^2+1"
$ws.Range("A10").Formula = "=3+(2+(3*2)*2+(2^2)--4)*+-+----3^2-4^2"
$ws.Range("A11").Formula = "=3+(2+(3*2)*2+(2^2)--4)*+-+----3"
$ws.Range("A12").Formula = "=3+(1)*+-+----3"
$ws.Range("A13").Formula = "=3+1*-3"
$ws.Range("A14").Formula = "=+--+((1*+-7)^2)"
$ws.Range("A15").Formula = "=+--+((1*+-7)^---2)"
$ws.Range("A16").Formula = "=--(+-(1+3-2/2)*(2^2^2)^2+1)*2^2-1+(1*(1-3))"

# G24 / G27 formulas are untouched in content, just restore them (their
# fill/style is refreshed below).
$ws.Range("G24").Formula = "=(2+2)*(2*2)^2"
$ws.Range("G27").Formula = "=1*(2*2)^2"

# ---------------------------------------------------------------------------
# 3. Formatting
#    - E3 keeps the original light-green highlight.
#    - A1:A16 (minus A9/A15 which get extra number formats) + G24 + G27 get
#      the new dark-green highlight.
#    - A9 additionally gets the accounting/comma number format (applied via
#      the built-in "Comma" style, same as the original workbook used).
#    - A15 additionally gets a 6-decimal-place number format.
# ---------------------------------------------------------------------------
$ws.Range("E3").Interior.Color = $lightGreen

$greenCells = "A1,A2,A3,A4,A5,A6,A7,A8,A9,A10,A11,A12,A13,A14,A15,A16,G24,G27"
foreach ($addr in $greenCells.Split(",")) {
    $ws.Range($addr).Interior.Color = $darkGreen
}

$ws.Range("A9").Style = "Comma"
$ws.Range("A9").Interior.Color = $darkGreen

$ws.Range("A15").NumberFormat = "0.000000"

# ---------------------------------------------------------------------------
# 4. Column A width (widened to fit the longer formula text).
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 31.5

# ---------------------------------------------------------------------------
# 5. Selection moves to A16, the newly-added final test case.
# ---------------------------------------------------------------------------
$ws.Range("A16").Select()

# ---------------------------------------------------------------------------
# Sheet2!A2's style index shifts automatically as a by-product of the style
# table changes above - nothing further required there.
# ---------------------------------------------------------------------------
